# Update column F (dSF) values to reflect repulled data / recalculated means.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "F2"  = -4
    "F3"  = 0
    "F4"  = 2
    "F5"  = -7
    "F6"  = -4
    "F7"  = -4
    "F8"  = -2
    "F9"  = -2
    "F10" = -5
    "F12" = -1
    "F13" = -1
    "F14" = 0
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
